$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the used range extent (header row 1, data rows 2..529, column O = timestamp)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = $ws.UsedRange.Rows.Count }

$newTimestamp = "2023-01-01 20:51:23"

# Update every timestamp cell in column O (rows 2..last) from the old crawl time to the new one.
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 15)
    if ($cell.Value2 -eq "2023-01-01 12:56:55") {
        $cell.Value = $newTimestamp
    }
}

# Row 320: product went out of stock online — annotate the aria-label text.
$ws.Cells.Item(320, 13).Value = "Fairtrade Papaya 1 Stück - Online kein Bestand 2.95 Schweizer Franken"
